$wb = $excel.ActiveWorkbook

# --- Clients sheet ---
$clients = $wb.Worksheets.Item("Clients")
$clients.Range("H3").Value = 1102
$clients.Range("I3").Value = 1202
$clients.Range("E4").Value = 2

# --- Episodes sheet ---
$episodes = $wb.Worksheets.Item("Episodes")
$episodes.Range("H3").Value = 2
# Drop the stray AA3 cell (no header in row 2); this also shrinks the
# sheet's used range / dimension back to A1:Z3 and row spans to 1:26.
$episodes.Range("AA3").ClearContents()

# --- K5 sheet ---
$k5 = $wb.Worksheets.Item("K5")
$k5.Range("E3").Value = 2

# --- Practitioners sheet ---
$practitioners = $wb.Worksheets.Item("Practitioners")
$practitioners.Range("C3").Value = 8
